$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2665998480081517
$ws.Range("D2").Value = 0.2665998480081517
$ws.Range("E2").Value = 0.29757278131017006
$ws.Range("F2").Value = 0.0019391581333836599
$ws.Range("G2").Value = 0.5997
$ws.Range("C3").Value = 8.04993745311914
$ws.Range("D3").Value = 8.04993745311914
$ws.Range("E3").Value = 8.985159951120172
$ws.Range("F3").Value = 0.05855255283179465
$ws.Range("G3").Value = 0.0016
$ws.Range("C4").Value = 0.15399949709097172
$ws.Range("D4").Value = 0.15399949709097172
$ws.Range("E4").Value = 0.17189079068164623
$ws.Range("F4").Value = 0.0011201408386092555
$ws.Range("G4").Value = 0.6833
$ws.Range("C5").Value = 129.01172595203948
$ws.Range("D5").Value = 0.8959147635558298
$ws.Range("F5").Value = 0.9383881481962124
$ws.Range("C6").Value = 137.48226275025775
$ws.Range("C7").Value = 2.6843161657517656
$ws.Range("D7").Value = 2.6843161657517656
$ws.Range("E7").Value = 0.869989211860278
$ws.Range("F7").Value = 0.005992891628625375
$ws.Range("G7").Value = 0.3491
$ws.Range("C8").Value = 0.10534124106987797
$ws.Range("D8").Value = 0.10534124106987797
$ws.Range("E8").Value = 0.03414118815959241
$ws.Range("F8").Value = 0.0002351804343360122
$ws.Range("G8").Value = 0.8546
$ws.Range("C9").Value = 0.8209136992146724
$ws.Range("D9").Value = 0.8209136992146724
$ws.Range("E9").Value = 0.26605884630772036
$ws.Range("F9").Value = 0.0018327374765370495
$ws.Range("G9").Value = 0.6112
$ws.Range("C10").Value = 444.3061162123164
$ws.Range("D10").Value = 3.0854591403633083
$ws.Range("F10").Value = 0.9919391904605015
$ws.Range("C11").Value = 447.91668731835273
$ws.Range("C12").Value = 0.2787109965819825
$ws.Range("D12").Value = 0.2787109965819825
$ws.Range("E12").Value = 0.3899439556233056
$ws.Range("F12").Value = 0.002674055261177791
$ws.Range("G12").Value = 0.5332
$ws.Range("C13").Value = 0.003316528784480466
$ws.Range("D13").Value = 0.003316528784480466
$ws.Range("E13").Value = 0.004640148286285702
$ws.Range("F13").Value = 0.000031819990433635044
$ws.Range("G13").Value = 0.9436
$ws.Range("C14").Value = 1.0223436065839182
$ws.Range("D14").Value = 1.0223436065839182
$ws.Range("E14").Value = 1.4303587402238187
$ws.Range("F14").Value = 0.009808738562323137
$ws.Range("G14").Value = 0.2384
$ws.Range("C15").Value = 102.92346612644043
$ws.Range("D15").Value = 0.7147462925447252
$ws.Range("F15").Value = 0.9874853861860655
$ws.Range("C16").Value = 104.2278372583908
$ws.Range("C17").Value = 0.02078980682888689
$ws.Range("D17").Value = 0.02078980682888689
$ws.Range("E17").Value = 0.07813472959791257
$ws.Range("F17").Value = 0.0005388276286779248
$ws.Range("G17").Value = 0.7736
$ws.Range("C18").Value = 0.24068105208406668
$ws.Range("D18").Value = 0.24068105208406668
$ws.Range("E18").Value = 0.9045562125089034
$ws.Range("F18").Value = 0.0062379415850065045
$ws.Range("G18").Value = 0.3506
$ws.Range("C19").Value = 0.006937995123558305
$ws.Range("D19").Value = 0.006937995123558305
$ws.Range("E19").Value = 0.02607520009169266
$ws.Range("F19").Value = 0.00017981809503932192
$ws.Range("G19").Value = 0.8685
$ws.Range("C20").Value = 38.31500024081087
$ws.Range("D20").Value = 0.2660763905611866
$ws.Range("F20").Value = 0.9930434126912762
$ws.Range("C21").Value = 38.58340909484738
$ws.Range("C22").Value = 0.7369877416089986
$ws.Range("D22").Value = 0.7369877416089986
$ws.Range("E22").Value = 2.0768337315805487
$ws.Range("F22").Value = 0.01382020410718638
$ws.Range("G22").Value = 0.1505
$ws.Range("C23").Value = 0.9884033681827945
$ws.Range("D23").Value = 0.9884033681827945
$ws.Range("E23").Value = 2.7853237436056593
$ws.Range("F23").Value = 0.01853482156798729
$ws.Range("G23").Value = 0.0984
$ws.Range("C24").Value = 0.5014287321337436
$ws.Range("D24").Value = 0.5014287321337436
$ws.Range("E24").Value = 1.4130277155023863
$ws.Range("F24").Value = 0.009402934447955287
$ws.Range("G24").Value = 0.2379
$ws.Range("C25").Value = 51.10001497853646
$ws.Range("D25").Value = 0.35486121512872537
$ws.Range("F25").Value = 0.958242039876871
$ws.Range("C26").Value = 53.326834820462
$ws.Range("C27").Value = 0.737121511836743
$ws.Range("D27").Value = 0.737121511836743
$ws.Range("E27").Value = 1.289584052386715
$ws.Range("F27").Value = 0.008853842695269116
$ws.Range("G27").Value = 0.2559
$ws.Range("C28").Value = 0.12136072294308964
$ws.Range("D28").Value = 0.12136072294308964
$ws.Range("E28").Value = 0.21231893301221882
$ws.Range("F28").Value = 0.0014577091199588195
$ws.Range("G28").Value = 0.6403
$ws.Range("C29").Value = 0.08606321918091919
$ws.Range("D29").Value = 0.08606321918091919
$ws.Range("E29").Value = 0.1505664305959868
$ws.Range("F29").Value = 0.001033737575474654
$ws.Range("G29").Value = 0.7039
$ws.Range("C30").Value = 82.30987154969914
$ws.Range("D30").Value = 0.571596330206244
$ws.Range("F30").Value = 0.9886547106092973
$ws.Range("C31").Value = 83.2544170036599
$ws.Range("C32").Value = 0.00965780160988203
$ws.Range("D32").Value = 0.00965780160988203
$ws.Range("E32").Value = 0.024196518914762974
$ws.Range("F32").Value = 0.00015243169510279552
$ws.Range("G32").Value = 0.8765
$ws.Range("C33").Value = 5.872276968323881
$ws.Range("D33").Value = 5.872276968323881
$ws.Range("E33").Value = 14.71231927060794
$ws.Range("F33").Value = 0.0926837357560556
$ws.Range("C34").Value = 0.00011190948327222117
$ws.Range("D34").Value = 0.00011190948327222117
$ws.Range("E34").Value = 0.0002803764291416964
$ws.Range("F34").Value = 0.0000017662976443632895
$ws.Range("G34").Value = 0.9863
$ws.Range("C35").Value = 57.47617815282069
$ws.Range("D35").Value = 0.39914012606125476
$ws.Range("F35").Value = 0.9071620662511972
$ws.Range("C36").Value = 63.358224832237724
$ws.Range("C37").Value = 0.023609627504653095
$ws.Range("D37").Value = 0.023609627504653095
$ws.Range("E37").Value = 0.41578310384559447
$ws.Range("F37").Value = 0.002777696038530913
$ws.Range("G37").Value = 0.3849
$ws.Range("C38").Value = 0.2769127307346845
$ws.Range("D38").Value = 0.2769127307346845
$ws.Range("E38").Value = 4.876639186981454
$ws.Range("F38").Value = 0.0325790567864282
$ws.Range("G38").Value = 0.0405
$ws.Range("C39").Value = 0.02236701553072412
$ws.Range("D39").Value = 0.02236701553072412
$ws.Range("E39").Value = 0.3938997825905656
$ws.Range("F39").Value = 0.002631501510187197
$ws.Range("G39").Value = 0.4431
$ws.Range("C40").Value = 8.176826641643892
$ws.Range("D40").Value = 0.05678351834474925
$ws.Range("F40").Value = 0.9620117456648537
$ws.Range("C41").Value = 8.499716015413954
